# Scheduled runner update: refresh currentAveragePrice / profit columns
# (H, I, J, K, L, M, N) for the Leve profit tables on each job sheet,
# per the latest market data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7144132
$ws.Range("I137").Value = 972.3333
$ws.Range("J137").Value = 20001820
$ws.Range("K137").Value = 2916.9999
$ws.Range("L137").Value = 60005460
$ws.Range("M137").Value = -366.9998999999998
$ws.Range("N137").Value = -60010560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2067575.4
$ws.Range("I45").Value = 3031764.5
$ws.Range("J45").Value = 1455.7142
$ws.Range("K45").Value = 3031764.5
$ws.Range("L45").Value = 1455.7142
$ws.Range("M45").Value = -3031387.5
$ws.Range("N45").Value = -2209.7142
$ws.Range("H61").Value = 25002306
$ws.Range("I61").Value = 31252226
$ws.Range("J61").Value = 2625
$ws.Range("K61").Value = 31252226
$ws.Range("L61").Value = 2625
$ws.Range("M61").Value = -31252014
$ws.Range("N61").Value = -3049
$ws.Range("H74").Value = 25005060
$ws.Range("I74").Value = 38464416
$ws.Range("K74").Value = 38464416
$ws.Range("M74").Value = -38463542
$ws.Range("H77").Value = 25005060
$ws.Range("I77").Value = 38464416
$ws.Range("K77").Value = 192322080
$ws.Range("M77").Value = -192317712
$ws.Range("H97").Value = 5182.273
$ws.Range("I97").Value = 6880.5
$ws.Range("J97").Value = 653.6667
$ws.Range("K97").Value = 6880.5
$ws.Range("L97").Value = 653.6667
$ws.Range("M97").Value = -6384.5
$ws.Range("N97").Value = -1645.6667
$ws.Range("H132").Value = 12502445
$ws.Range("I132").Value = 17859328
$ws.Range("J132").Value = 3052
$ws.Range("K132").Value = 53577984
$ws.Range("L132").Value = 9156
$ws.Range("M132").Value = -53575454
$ws.Range("N132").Value = -14216
$ws.Range("H136").Value = 25002306
$ws.Range("I136").Value = 31252226
$ws.Range("J136").Value = 2625
$ws.Range("K136").Value = 93756678
$ws.Range("L136").Value = 7875
$ws.Range("M136").Value = -93754128
$ws.Range("N136").Value = -12975

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 842.8570999999999
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -350
$ws.Range("N22").Value = -1600
$ws.Range("H31").Value = 7095837
$ws.Range("I31").Value = 3697.375
$ws.Range("J31").Value = 47622348
$ws.Range("K31").Value = 3697.375
$ws.Range("L31").Value = 47622348
$ws.Range("M31").Value = -3402.375
$ws.Range("N31").Value = -47622938
$ws.Range("H34").Value = 7095837
$ws.Range("I34").Value = 3697.375
$ws.Range("J34").Value = 47622348
$ws.Range("K34").Value = 3697.375
$ws.Range("L34").Value = 47622348
$ws.Range("M34").Value = -3495.375
$ws.Range("N34").Value = -47622752
$ws.Range("H62").Value = 2357.8262
$ws.Range("I62").Value = 2178.2354
$ws.Range("J62").Value = 2866.6667
$ws.Range("K62").Value = 2178.2354
$ws.Range("L62").Value = 2866.6667
$ws.Range("M62").Value = -1554.2354
$ws.Range("N62").Value = -4114.6667
$ws.Range("H65").Value = 2357.8262
$ws.Range("I65").Value = 2178.2354
$ws.Range("J65").Value = 2866.6667
$ws.Range("K65").Value = 10891.177
$ws.Range("L65").Value = 14333.3335
$ws.Range("M65").Value = -7771.177
$ws.Range("N65").Value = -20573.3335
$ws.Range("H134").Value = 1703602.8
$ws.Range("I134").Value = 2990.4443
$ws.Range("K134").Value = 8971.332900000001
$ws.Range("M134").Value = -6436.332900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1840.3334
$ws.Range("I5").Value = 336.83334
$ws.Range("K5").Value = 1010.50002
$ws.Range("M5").Value = -898.5000200000001
$ws.Range("H132").Value = 1351.4736
$ws.Range("I132").Value = 772.25
$ws.Range("J132").Value = 1772.7273
$ws.Range("K132").Value = 6950.25
$ws.Range("L132").Value = 15954.5457
$ws.Range("M132").Value = -4420.25
$ws.Range("N132").Value = -21014.5457
$ws.Range("H134").Value = 4089.9688
$ws.Range("I134").Value = 1993.2778
$ws.Range("J134").Value = 6785.7144
$ws.Range("K134").Value = 5979.8334
$ws.Range("L134").Value = 20357.1432
$ws.Range("M134").Value = -909.8334000000004
$ws.Range("N134").Value = -30497.1432
$ws.Range("H135").Value = 1840.3334
$ws.Range("I135").Value = 336.83334
$ws.Range("K135").Value = 3031.50006
$ws.Range("M135").Value = -496.5000600000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15110
$ws.Range("I70").Value = 37166.668
$ws.Range("J70").Value = 4930
$ws.Range("K70").Value = 37166.668
$ws.Range("L70").Value = 4930
$ws.Range("M70").Value = -36896.668
$ws.Range("N70").Value = -5470
$ws.Range("H73").Value = 15110
$ws.Range("I73").Value = 37166.668
$ws.Range("J73").Value = 4930
$ws.Range("K73").Value = 37166.668
$ws.Range("L73").Value = 4930
$ws.Range("M73").Value = -36230.668
$ws.Range("N73").Value = -6802
$ws.Range("H113").Value = 77687.38
$ws.Range("I113").Value = 111847
$ws.Range("J113").Value = 828.25
$ws.Range("K113").Value = 111847
$ws.Range("L113").Value = 828.25
$ws.Range("M113").Value = -109677
$ws.Range("N113").Value = -5168.25
$ws.Range("H122").Value = 2779521
$ws.Range("I122").Value = 3922769
$ws.Range("J122").Value = 3062
$ws.Range("K122").Value = 11768307
$ws.Range("L122").Value = 9186
$ws.Range("M122").Value = -11765857
$ws.Range("N122").Value = -14086
$ws.Range("H132").Value = 4820.8423
$ws.Range("I132").Value = 4707.8184
$ws.Range("K132").Value = 14123.4552
$ws.Range("M132").Value = -11593.4552
$ws.Range("H133").Value = 61374.938
$ws.Range("J133").Value = 61374.938
$ws.Range("L133").Value = 61374.938
$ws.Range("N133").Value = -71494.93799999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1050.375
$ws.Range("I93").Value = 950.5
$ws.Range("J93").Value = 1350
$ws.Range("K93").Value = 950.5
$ws.Range("L93").Value = 1350
$ws.Range("M93").Value = 297.5
$ws.Range("N93").Value = -3846
$ws.Range("H121").Value = 50420
$ws.Range("J121").Value = 50420
$ws.Range("L121").Value = 50420
$ws.Range("N121").Value = -53914

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2437.1365
$ws.Range("I126").Value = 1827.4667
$ws.Range("J126").Value = 3743.5715
$ws.Range("K126").Value = 5482.4001
$ws.Range("L126").Value = 11230.7145
$ws.Range("M126").Value = -3012.4001
$ws.Range("N126").Value = -16170.7145
